$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.927.56"
Set-TextValue "E2" "  -1.04%  "
Set-TextValue "D3" "1.808.10"
Set-TextValue "E3" "  -0.61%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.19%  "
Set-TextValue "D5" "310.34"
Set-TextValue "E5" "  -0.92%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  -0.09%  "
Set-TextValue "D7" "0.4630"
Set-TextValue "E7" "  +3.53%  "
Set-TextValue "E8" "  -1.67%  "
Set-TextValue "D9" "0.07372"
Set-TextValue "E9" "  -0.23%  "
Set-TextValue "D10" "0.8742"
Set-TextValue "E10" "  -0.61%  "
Set-TextValue "E11" "  -1.83%  "
Set-TextValue "D12" "1.810.65"
Set-TextValue "E12" "  -0.47%  "
Set-TextValue "D13" "5.370"
Set-TextValue "E13" "  -0.86%  "
Set-TextValue "D14" "92.59"
Set-TextValue "E14" "  -0.37%  "
Set-TextValue "D15" "6.512"
Set-TextValue "E15" "  -2.91%  "
Set-TextValue "D16" "0.07030"
Set-TextValue "E16" "  -1.34%  "
Set-TextValue "E17" "  -0.25%  "
Set-TextValue "E18" "  -0.79%  "
Set-TextValue "E19" "  -0.11%  "
Set-TextValue "D20" "14.71"
Set-TextValue "E20" "  -2.19%  "
Set-TextValue "D21" "26.938.63"
Set-TextValue "E21" "  -1.09%  "
Set-TextValue "D22" "5.297"
Set-TextValue "E22" "  -1.12%  "
Set-TextValue "D23" "10.63"
Set-TextValue "E23" "  -2.62%  "
Set-TextValue "D24" "2.046.77"
Set-TextValue "E24" "  -0.34%  "
Set-TextValue "D25" "1.910"
Set-TextValue "E25" "  -2.90%  "
Set-TextValue "D26" "151.62"
Set-TextValue "E26" "  +0.27%  "
Set-TextValue "D27" "18.33"
Set-TextValue "E27" "  -1.15%  "
Set-TextValue "D28" "2.141"
Set-TextValue "E28" "  -6.43%  "
Set-TextValue "D29" "5.289"
Set-TextValue "E29" "  -0.97%  "
Set-TextValue "D30" "115.88"
Set-TextValue "E30" "  -1.13%  "
Set-TextValue "D31" "0.08915"
Set-TextValue "D32" "0.7557"
Set-TextValue "E32" "  -3.17%  "
Set-TextValue "E33" "  -2.80%  "
Set-TextValue "D34" "2.929"
Set-TextValue "E34" "  +0.69%  "
Set-TextValue "D35" "4.462"
Set-TextValue "E35" "  -2.20%  "
Set-TextValue "D36" "1.000"
Set-TextValue "E36" "  -0.14%  "
Set-TextValue "D37" "1.105"
Set-TextValue "E37" "  -0.56%  "
Set-TextValue "D38" "0.01968"
Set-TextValue "E38" "  -0.45%  "
Set-TextValue "D39" "0.05248"
Set-TextValue "E39" "  -0.10%  "
Set-TextValue "D40" "2.422"
Set-TextValue "E40" "  +6.33%  "
Set-TextValue "D41" "2.925"
Set-TextValue "E41" "  +2.04%  "
Set-TextValue "D42" "7.215"
Set-TextValue "E42" "  -1.39%  "
Set-TextValue "D43" "0.5310"
Set-TextValue "E43" "  +0.44%  "
Set-TextValue "D44" "0.1665"
Set-TextValue "E44" "  -2.23%  "
Set-TextValue "D45" "8.519"
Set-TextValue "E45" "  -0.65%  "
Set-TextValue "D46" "0.5000"
Set-TextValue "E46" "  -0.73%  "
Set-TextValue "D47" "10.39"
Set-TextValue "E47" "  -1.42%  "
Set-TextValue "D48" "104.01"
Set-TextValue "E48" "  -0.77%  "
Set-TextValue "B49" "PaxDollar"
Set-TextValue "C49" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D49" "0.9999"
Set-TextValue "E49" "  -0.11%  "
Set-TextValue "B50" "NEARProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D50" "1.669"
Set-TextValue "E50" "  -0.90%  "
Set-TextValue "E51" "  -1.60%  "
